$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.668.96"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").Value = "2.119.23"
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.015"
$ws.Range("E4").Value = "  +1.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "338.77"
$ws.Range("E5").Value = "  +1.42%  "
$ws.Range("E6").Value = "  +1.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5276"
$ws.Range("E7").Value = "  +0.79%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4552"
$ws.Range("E8").Value = "  +1.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.04"
$ws.Range("E9").Value = "  +1.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09113"
$ws.Range("E10").Value = "  +1.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.174"
$ws.Range("E11").Value = "  +0.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.47"
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("D13").Value = "2.120.80"
$ws.Range("E13").Value = "  +0.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.839"
$ws.Range("E14").Value = "  +0.75%  "
$ws.Range("E15").Value = "  +3.98%  "
$ws.Range("E16").Value = "  +2.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001171"
$ws.Range("E17").Value = "  +3.97%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.014"
$ws.Range("E18").Value = "  +1.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06703"
$ws.Range("E19").Value = "  +1.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.55"
$ws.Range("E20").Value = "  +1.18%  "
$ws.Range("E21").Value = "  +1.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.452"
$ws.Range("E22").Value = "  +2.28%  "
$ws.Range("D23").Value = "30.732.16"
$ws.Range("E23").Value = "  +0.62%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.93"
$ws.Range("E24").Value = "  +4.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.379"
$ws.Range("E25").Value = "  +1.54%  "
$ws.Range("D26").Value = "2.376.37"
$ws.Range("E26").Value = "  +0.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.51"
$ws.Range("E27").Value = "  +0.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "165.63"
$ws.Range("E28").Value = "  +1.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.551"
$ws.Range("E29").Value = "  -1.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "135.92"
$ws.Range("E30").Value = "  +2.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.207"
$ws.Range("E31").Value = "  +0.33%  "
$ws.Range("E32").Value = "  +0.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.419"
$ws.Range("E33").Value = "  +3.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.640"
$ws.Range("E34").Value = "  -1.80%  "
$ws.Range("E35").Value = "  +0.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.55"
$ws.Range("E36").Value = "  +0.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.961"
$ws.Range("E37").Value = "  +8.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02668"
$ws.Range("E38").Value = "  +3.66%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06880"
$ws.Range("E39").Value = "  +1.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2325"
$ws.Range("E40").Value = "  +2.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.63"
$ws.Range("E41").Value = "  -1.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6898"
$ws.Range("E42").Value = "  -0.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.266"
$ws.Range("E43").Value = "  +0.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.19"
$ws.Range("E44").Value = "  +7.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6484"
$ws.Range("E45").Value = "  +1.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.311"
$ws.Range("E46").Value = "  -1.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000367"
$ws.Range("E47").Value = "  +15.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.703"
$ws.Range("E48").Value = "  +1.42%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.260"
$ws.Range("E49").Value = "  +0.95%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "83.08"
$ws.Range("E50").Value = "  -0.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07312"
$ws.Range("E51").Value = "  +3.42%  "
